# Added page load validation (test will fail if > 25 sec)
# Appends the newly captured submission-time samples to the four
# "RF LOAD RESULTS" sheets.

$wb = $excel.ActiveWorkbook

function Add-Rows {
    param(
        [string]$SheetName,
        [int]$StartRow,
        [object[]]$Rows
    )

    $ws = $wb.Worksheets.Item($SheetName)
    $r = $StartRow
    foreach ($row in $Rows) {
        $ws.Cells.Item($r, 1).Value = $row[0]
        $ws.Cells.Item($r, 2).Value = $row[1]
        $ws.Cells.Item($r, 3).Value = $row[2]
        $ws.Cells.Item($r, 4).Value = $row[3]
        $ws.Cells.Item($r, 5).Value = $row[4]
        $r = $r + 1
    }
}


# --- Sheet 1: "Submit orders" (rows 66-69) ---
$rows1 = @(
    ,@("09.19.2022 16:17 (Kyiv+Israel) 13:17 (UTC) 22:17 (Japan) 18:47 (India)", 0.912, -0.09900000000000009, "***", "***")
    ,@("09.21.2022 09:39 (Kyiv+Israel) 06:39 (UTC) 15:39 (Japan) 12:09 (India)", "***", "***", 1.644, -0.506)
    ,@("09.26.2022 09:05 (Kyiv+Israel) 06:05 (UTC) 15:05 (Japan) 11:35 (India)", "***", "***", 1.2, -0.06200000000000006)
    ,@("09.27.2022 09:33 (Kyiv+Israel) 06:33 (UTC) 15:33 (Japan) 12:03 (India)", 2.18, -1.367, "***", "***")
)
Add-Rows "Submit orders" 66 $rows1

# --- Sheet 2: "Submit internet survey" (rows 66-69) ---
$rows2 = @(
    ,@("09.19.2022 16:19 (Kyiv+Israel) 13:19 (UTC) 22:19 (Japan) 18:49 (India)", 1.022, -0.355, "***", "***")
    ,@("09.21.2022 09:42 (Kyiv+Israel) 06:42 (UTC) 15:42 (Japan) 12:12 (India)", "***", "***", 0.9370000000000001, -0.01000000000000001)
    ,@("09.26.2022 09:09 (Kyiv+Israel) 06:09 (UTC) 15:09 (Japan) 11:39 (India)", "***", "***", 0.765, 0.162)
    ,@("09.27.2022 09:35 (Kyiv+Israel) 06:35 (UTC) 15:35 (Japan) 12:05 (India)", 0.631, 0.03600000000000003, "***", "***")
)
Add-Rows "Submit internet survey" 66 $rows2

# --- Sheet 3: "Submit a phone survey" (rows 59-62) ---
$rows3 = @(
    ,@("09.19.2022 16:22 (Kyiv+Israel) 13:22 (UTC) 22:22 (Japan) 18:52 (India)", 2.754, -1.29, "***", "***")
    ,@("09.21.2022 09:44 (Kyiv+Israel) 06:44 (UTC) 15:44 (Japan) 12:14 (India)", "***", "***", 2.011, -0.2070000000000001)
    ,@("09.26.2022 09:12 (Kyiv+Israel) 06:12 (UTC) 15:12 (Japan) 11:42 (India)", "***", "***", 1.966, -0.1619999999999999)
    ,@("09.27.2022 09:38 (Kyiv+Israel) 06:38 (UTC) 15:38 (Japan) 12:08 (India)", 2.656, -1.192, "***", "***")
)
Add-Rows "Submit a phone survey" 59 $rows3

# --- Sheet 4: "Checkertificate" (rows 56-64) ---
$rows4 = @(
    ,@("09.19.2022 16:30 (Kyiv+Israel) 13:30 (UTC) 22:30 (Japan) 19:00 (India)", 0.828, -0.102, "***", "***")
    ,@("09.21.2022 09:32 (Kyiv+Israel) 06:32 (UTC) 15:32 (Japan) 12:02 (India)", 1.243, -0.5170000000000001, "***", "***")
    ,@("09.21.2022 09:47 (Kyiv+Israel) 06:47 (UTC) 15:47 (Japan) 12:17 (India)", "***", "***", 1.256, -0.334)
    ,@("09.22.2022 18:26 (Kyiv+Israel) 15:26 (UTC) 00:26 (Japan) 20:56 (India)", 0.873, -0.147, "***", "***")
    ,@("09.24.2022 20:55 (Kyiv+Israel) 17:55 (UTC) 02:55 (Japan) 23:25 (India)", 0.8100000000000001, -0.08400000000000007, "***", "***")
    ,@("09.26.2022 09:15 (Kyiv+Israel) 06:15 (UTC) 15:15 (Japan) 11:45 (India)", "***", "***", 1.096, -0.174)
    ,@("09.26.2022 09:52 (Kyiv+Israel) 06:52 (UTC) 15:52 (Japan) 12:22 (India)", 0.869, -0.143, "***", "***")
    ,@("09.26.2022 17:11 (Kyiv+Israel) 14:11 (UTC) 23:11 (Japan) 19:41 (India)", 0.747, -0.02100000000000002, "***", "***")
    ,@("09.27.2022 09:44 (Kyiv+Israel) 06:44 (UTC) 15:44 (Japan) 12:14 (India)", 0.8139999999999999, -0.08799999999999997, "***", "***")
)
Add-Rows "Checkertificate" 56 $rows4
